$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2663136666666667
$ws.Range("N2").Value = 0.798941
$ws.Range("O2").Value = 0.2749240117919797
$ws.Range("P2").Value = 0.2749240117919797
$ws.Range("Q2").Value = 17.02042130819189
$ws.Range("R2").Value = 153.183791773727
$ws.Range("S2").Value = 0.1118370785892944
$ws.Range("T2").Value = 0.1118370785892944

$ws.Range("G3").Value = 63.91118233333333
$ws.Range("H3").Value = 191.733547
$ws.Range("I3").Value = 0.4067926910433548
$ws.Range("J3").Value = 0.4067926910433549
$ws.Range("M3").Value = 0.7023673333333335
$ws.Range("N3").Value = 2.107102
$ws.Range("O3").Value = 0.7250759882080203
$ws.Range("P3").Value = 0.7250759882080203
$ws.Range("Q3").Value = 44.88912670564378
$ws.Range("R3").Value = 404.002140350794
$ws.Range("S3").Value = 0.2949556124540604
$ws.Range("T3").Value = 0.2949556124540604

$ws.Range("I4").Value = 0.3656254573230189
$ws.Range("J4").Value = 0.365625457323019
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.2663136666666667
$ws.Range("N4").Value = 0.798941
$ws.Range("O4").Value = 0.2749240117919797
$ws.Range("P4").Value = 0.2749240117919797
$ws.Range("Q4").Value = 15.2979624798
$ws.Range("R4").Value = 137.6816623182
$ws.Range("S4").Value = 0.1005192175405216
$ws.Range("T4").Value = 0.1005192175405216

$ws.Range("I5").Value = 0.3656254573230189
$ws.Range("J5").Value = 0.365625457323019
$ws.Range("M5").Value = 0.7023673333333335
$ws.Range("N5").Value = 2.107102
$ws.Range("O5").Value = 0.7250759882080203
$ws.Range("P5").Value = 0.7250759882080203
$ws.Range("Q5").Value = 40.34636767560001
$ws.Range("R5").Value = 363.1173090804
$ws.Range("S5").Value = 0.2651062397824973
$ws.Range("T5").Value = 0.2651062397824974

$ws.Range("G6").Value = 35.755375
$ws.Range("H6").Value = 107.266125
$ws.Range("I6").Value = 0.2275818516336261
$ws.Range("J6").Value = 0.2275818516336262
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.2663136666666667
$ws.Range("N6").Value = 0.798941
$ws.Range("O6").Value = 0.2749240117919797
$ws.Range("P6").Value = 0.2749240117919797
$ws.Range("Q6").Value = 9.522145019291667
$ws.Range("R6").Value = 85.699305173625
$ws.Range("S6").Value = 0.0625677156621636
$ws.Range("T6").Value = 0.06256771566216361

$ws.Range("G7").Value = 35.755375
$ws.Range("H7").Value = 107.266125
$ws.Range("I7").Value = 0.2275818516336261
$ws.Range("J7").Value = 0.2275818516336262
$ws.Range("M7").Value = 0.7023673333333335
$ws.Range("N7").Value = 2.107102
$ws.Range("O7").Value = 0.7250759882080203
$ws.Range("P7").Value = 0.7250759882080203
$ws.Range("Q7").Value = 25.11340739108334
$ws.Range("R7").Value = 226.02066651975
$ws.Range("S7").Value = 0.1650141359714625
$ws.Range("T7").Value = 0.1650141359714626

